$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Update VENTA (D) values for rows 2 and 3
$ws.Range("D2").Value = 2626.18
$ws.Range("D3").Value = 821.23

# Recompute dependent values: POR CUMPLIR (E) = PRESUPUESTO (C) - VENTA (D)
$ws.Range("E2").Value = -2626.18
$ws.Range("E3").Value = 12902.11

# Recompute CUMPLIMIENTO (F) = VENTA (D) / PRESUPUESTO (C)
$ws.Range("F3").Value = 0.05984184608120181

# Update TOTAL row (row 4), which sums rows 2 and 3
$ws.Range("D4").Value = 3447.41
$ws.Range("E4").Value = 10275.93
$ws.Range("F4").Value = 0.2512077963527829
